$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Give row 18 (which will become the last/only data row of the table)
#    the "closing" row style currently used by row 22 (thicker bottom border).
$ws.Range("B22:J22").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# 2) Turn rows 16 and 17 into the two periods for employee
#    MAURICIO ANDRES FERNANDEZ HERRERA (1047499869), replacing the old
#    CARLOS ALFREDO CURE CAMARGO period rows.
$ws.Range("C16").Value = "1047499869"
$ws.Range("D16").Value = "MAURICIO ANDRES FERNANDEZ HERRERA"
$ws.Range("E16").Value = "1903"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 877803

$ws.Range("C17").Value = "1047499869"
$ws.Range("D17").Value = "MAURICIO ANDRES FERNANDEZ HERRERA"
$ws.Range("E17").Value = "1904"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 877803

# 3) Remove the now-obsolete rows (old CC/2502, FRANKLIN, and the old
#    MAURICIO rows that lived further down) - row 18 (CARLOS/2503) is kept
#    as the final table row.
$ws.Rows("19:22").Delete()

# 4) Update the summary figures at the top of the statement.
$ws.Range("E11").Value = 123190
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3
